# Updated cryptos list on Wed Jun 28 15:49:59 UTC 2023 with GitHub Actions
#
# Applies the per-cell Price (D) / Volume(1h) (E) refresh captured in the
# source diff, including the BitcoinCash <-> Avalanche row swap (rows 18/19:
# name, link, price and volume all move together).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.409.28"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.860.16"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.93"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06443"
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").Value = "1.881.53"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07436"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.35"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.009"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.70"
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6353"
$ws.Range("E15").Value = "  -3.91%  "
$ws.Range("D16").Value = "30.364.68"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9994"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.82"
$ws.Range("E18").Value = "  -2.99%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.08"
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007419"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("D21").Value = "2.100.55"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.019"
$ws.Range("E23").Value = "  -4.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.024"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.246"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.95"
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("E27").Value = "  -2.44%  "
$ws.Range("E28").Value = "  -1.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1040"
$ws.Range("E29").Value = "  +7.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.394"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.159"
$ws.Range("E31").Value = "  -3.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.940"
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04916"
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.155"
$ws.Range("E34").Value = "  -4.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7281"
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9998"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.695"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01898"
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.651"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9135"
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.975"
$ws.Range("E41").Value = "  -4.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.09"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9997"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4122"
$ws.Range("E44").Value = "  -3.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.586"
$ws.Range("E45").Value = "  -2.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.157"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "61.23"
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1214"
$ws.Range("E48").Value = "  -5.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.686"
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.412"
$ws.Range("E50").Value = "  -3.21%  "
$ws.Range("E51").Value = "  -0.37%  "
